$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps in column F (rows 2-15) on the data sheet ---
$newTimes = @(
    "2021-10-05 14:33:14.315643",
    "2021-10-05 14:33:14.315651",
    "2021-10-05 14:33:14.315654",
    "2021-10-05 14:33:14.315656",
    "2021-10-05 14:33:14.315659",
    "2021-10-05 14:33:14.315662",
    "2021-10-05 14:33:14.315664",
    "2021-10-05 14:33:14.315666",
    "2021-10-05 14:33:14.315669",
    "2021-10-05 14:33:14.315672",
    "2021-10-05 14:33:14.315674",
    "2021-10-05 14:33:14.315677",
    "2021-10-05 14:33:14.315679",
    "2021-10-05 14:33:14.315682"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet right after the "data" sheet ---
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1), columns B through G, styled like the "data" sheet headers
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Copy the header style (bold, bordered, centered) from the data sheet's header row
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2)
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "Atypical Haemolytic Uraemic Syndrome_MPGN"
$ws.Range("C2").Value = 211

# D2 holds the textual string "0.38" (not a number) - force text type, then
# strip the formatting override so the cell keeps the default (unstyled) look
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.38"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "2021-02-11T09:53:03.889058Z"
$ws.Range("F2").Value = "2021-10-05 14:33:14.311983"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/211/?format=json"

Write-Output "done"
